$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.067.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.241.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.38%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.76%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.240.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.459"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.34%  "

$ws.Range("E10").Value = "  -4.25%  "

$ws.Range("E11").Value = "  -5.62%  "

$ws.Range("E12").Value = "  -4.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.796.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.241.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.46%  "

$ws.Range("E17").Value = "  -6.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "59.101.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.41%  "

$ws.Range("E20").Value = "  -6.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "362.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.20%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.373.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.90%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0975"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.43%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.84%  "

$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.29%  "

$ws.Range("E39").Value = "  -6.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0709"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.273.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("E44").Value = "  -5.31%  "

$ws.Range("E45").Value = "  -3.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.49%  "

$ws.Range("E47").Value = "  -6.65%  "

$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.298.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.90%  "
